$wb = $excel.ActiveWorkbook

# Update the "想去人数" (want-to-go count) values for the first two
# convention rows on both the "展览" sheet and the "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5438
    $ws.Range("F3").Value = 166
}
